$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 / cell A10: correct the recalculated date-time serial (sub-second
# rounding noise from the WSL export run) from 45878.37517072909 to
# 45878.37517072917 — both represent 2025-08-09 09:00:14.751.
$ws.Range("A10").Value = 45878.37517072917

# New row 11: latest DropControl sensor reading for 2025-08-09 11:00:13.
$ws.Range("A11").Value = 45878.45849049932
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat
$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 37
$ws.Range("D11").Value = 15.89
$ws.Range("E11").Value = 85.83
$ws.Range("F11").Value = 671.48
$ws.Range("G11").Value = 16.77
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "11:00:13"
